$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Step 1: shift existing columns to make room for the two new columns ---
# Insert a new column at B (ProjectID) - pushes ProjectName..Progress right by one (B->C ... F->G)
$ws.Columns.Item(2).Insert()
# After that insert, Progress is now column G (7). Insert a new column at G
# (ProjectDependency) so Progress moves to H.
$ws.Columns.Item(7).Insert()

# --- Step 2: grow the table to cover the now-8-column range ---
$tbl.Resize($ws.Range("A1:H6"))

# --- Step 3: write the final header row (this also renames/retitles the table's
# ListColumns, since table column names are sourced from the header cells) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "ProjectID"
$ws.Range("C1").Value = "ProjectName"
$ws.Range("D1").Value = "TaskName"
$ws.Range("E1").Value = "EstimatedEffortHours"
$ws.Range("F1").Value = "TaskDependencies"
$ws.Range("G1").Value = "ProjectDependency"
$ws.Range("H1").Value = "Progress"

# --- Step 4: write the final data grid explicitly (row by row) ---
# Row 2: Project A / Design UI
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Project A"
$ws.Range("D2").Value = "Design UI"
$ws.Range("E2").Value = 100
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").Value = 50

# Row 3: Project A / Implement Backend
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Project A"
$ws.Range("D3").Value = "Implement Backend"
$ws.Range("E3").Value = 55
$ws.Range("F3").Value = 1
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()

# Row 4: Project A / Testing
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Project A"
$ws.Range("D4").Value = "Testing"
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 2
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# Row 5: Project B / Database Setup
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Project B"
$ws.Range("D5").Value = "Database Setup"
$ws.Range("E5").Value = 60
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 1
$ws.Range("H5").ClearContents()

# Row 6: Project B / API Development
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Project B"
$ws.Range("D6").Value = "API Development"
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 1
$ws.Range("H6").ClearContents()

# --- Step 5: cosmetic - match the author's final selection / column widths ---
# Columns 1,3,4,5,6,8 already carry over their original (exact) widths from the
# Insert() shift above, so only the two brand-new columns need an explicit width:
# column B (ProjectID) matches column A's width, column G (ProjectDependency)
# matches column F's width.
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 14.5

$ws.Range("G7").Select()
